$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Type" column (J) - a free-text field duplicating the "Rule For" /
# data-validation column - is being removed entirely from the import
# template. Its header cell comment ("Mandatory - use existing
# nomenclature...") must go with it; simply deleting the column leaves
# the comment behind (Excel re-homes it onto whatever shifts into J),
# so drop the comment explicitly first.
$ws.Range("J1").Comment.Delete()

# Delete the whole column; everything to the right (the "Rule For" /
# "Accounting" / "Reporting" column) shifts left to become the new
# column J, and the shared-string table / dimension / used range are
# recalculated accordingly.
$ws.Columns("J").Delete()
